{"js": "// Replace the \"[**give an example**]\" placeholder in the foldVector bullet\n// with the fleshed-out example text, then move the \"_GoBack\" bookmark from\n// its old position (just before \"toterror\" further down the document) to\n// sit right after the text we just edited -- mirroring Word's own behavior\n// of tracking the most recent edit location with that bookmark.\n\nconst body = context.document.body;\n\n// 1) Find and replace the placeholder text.\nconst placeholder = body.search(\"[**give an example**]\", { matchWildcards: false });\nplaceholder.load(\"items\");\nawait context.sync();\n\nif (placeholder.items.length === 0) {\n  throw new Error('Could not find placeholder text \"[**give an example**]\"');\n}\n\nconst target = placeholder.items[0];\ntarget.insertText(\n  \"For example, spatially autocorrelated regions could be treated as \\\"folds\\\" and thereby left out of the model fitting process (and used as an out-of-bag validation set) one at a time. This strategy can result in a more honest appraisal of predictive performance.\",\n  \"Replace\"\n);\nawait context.sync();\n\n// 2) Drop the old \"_GoBack\" bookmark (previously located right before the\n//    \"toterror\" run) -- a document can only have one bookmark of a given name.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3) Re-insert \"_GoBack\" as a collapsed bookmark immediately after the text\n//    we just inserted, so it now marks the latest edit location.\nconst afterEdit = target.getRange(\"End\");\nafterEdit.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Replace the \"[**give an example**]\" placeholder in the foldVector bullet\n# with the fleshed-out example text, then move the \"_GoBack\" bookmark from\n# its old position (just before \"toterror\" further down the document) to\n# sit right after the text we just edited -- mirroring Word's own behavior\n# of tracking the most recent edit location with that bookmark.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: find & replace the placeholder text -----------------------\n$range = $d.Content\n$found = $range.Find.Execute(\"[**give an example**]\")\nif (-not $found) {\n    throw 'Could not find placeholder text \"[**give an example**]\"'\n}\n\n# Pin the boundary between the preceding run and the placeholder run with a\n# throwaway bookmark spanning the whole placeholder BEFORE we touch the\n# text. Without this, replacing $range.Text merges the edited run into the\n# identically-formatted run right before it; with the pin in place, the\n# replacement text lands in its own run, matching Word's real behavior.\n$d.Bookmarks.Add(\"ZZZ_TMP_PIN\", $range)\n$range.Text = \"For example, spatially autocorrelated regions could be treated as `\"folds`\" and thereby left out of the model fitting process (and used as an out-of-bag validation set) one at a time. This strategy can result in a more honest appraisal of predictive performance.\"\n\n$pin = $d.Bookmarks.Item(\"ZZZ_TMP_PIN\")\n$afterEditPos = $pin.Range.End\n$pin.Delete()\n\n# --- Step 2: drop the old \"_GoBack\" bookmark ----------------------------\n# (it used to sit right before the \"toterror\" run, further down the doc;\n# a document can only have one bookmark of a given name at a time).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- Step 3: re-create \"_GoBack\" right after the text we just edited ----\n# Build it as a zero-length bookmark: insert a throwaway character right\n# after the new text, wrap a bookmark around that character (Bookmarks.Add\n# needs a non-empty range to anchor reliably), then clear the character via\n# the bookmark's own range so the bookmark collapses in place.\n$insPoint = $d.Range($afterEditPos, $afterEditPos)\n$insPoint.InsertAfter(\"#\")\n$markRange = $d.Range($afterEditPos, $afterEditPos + 1)\n$d.Bookmarks.Add(\"_GoBack\", $markRange)\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$bm.Range.Text = \"\"\n"}
